# Daily attendance processing - 2026-01-15 17:22:53
# Reorders the "Recorded By" (column G) entries so that when the last
# author in the comma-separated list is "System" (any case), that
# leading author is rotated to the end of the list, e.g.
#   "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, system, System" -> "system, System, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $rawParts = $text.Split(",")
    if ($rawParts.Count -le 1) {
        continue
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $lastPart = $parts[$parts.Count - 1]
    if ($lastPart.ToLower() -eq "system") {
        $first = $parts[0]
        $rest = $parts[1..($parts.Count - 1)]
        $newParts = $rest + ,$first
        $newValue = [string]::Join(", ", $newParts)
        $cell.Value = $newValue
    }
}
